$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update trigger time cells
$ws.Range("B3").Value = "17:05 - 17:09"
$ws.Range("B6").Value = "18:55 - 18:59"
$ws.Range("B7").Value = "19:00 - 19:04"
$ws.Range("B12").Value = "11:55 - 11:59"

# Update the message text (wrap with asterisks) and trigger JSON array
$ws.Range("A7").Value = "*Отче наш, Отец Небесный! Волею Создателя, Пророка и Народа Пространство Святая Русь ЕСМЬ Равенство и Любовь Навечно! Да будет Свет Истины!*"
$ws.Range("C12").Value = '["08", "13", "26"]'

# Update view state: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B20").Select()
